$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.057.32"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "3.529.34"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'588.33"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "'178.05"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.530.32"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").Value = "4.141.70"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "'30.55"
$ws.Range("E14").Value = "  -4.20%  "
$ws.Range("D16").Value = "67.070.24"
$ws.Range("D17").Value = "'0.0000174"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "3.536.87"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").Value = "'14.08"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'384.18"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "'7.87"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").Value = "'0.542"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D26").Value = "'72.03"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").Value = "'0.0000122"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").Value = "'9.98"
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'24.58"
$ws.Range("E31").Value = "  +4.60%  "
$ws.Range("D32").Value = "'5.97"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'29.53"
$ws.Range("E38").Value = "  +13.02%  "
$ws.Range("D39").Value = "'159.84"
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("E40").Value = "  +3.38%  "
$ws.Range("D41").Value = "'1.81"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'4.55"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.60"
$ws.Range("E44").Value = "  -5.49%  "
$ws.Range("D45").Value = "2.754.31"
$ws.Range("E45").Value = "  -2.77%  "
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").Value = "'25.57"
$ws.Range("E47").Value = "  -5.39%  "
$ws.Range("D48").Value = "'40.85"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "'326.38"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("E51").Value = "  -1.70%  "
